$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 841, shifting the existing rows 841:889 down
# to become 845:893 (dimension grows from R889 to R893).
$ws.Rows("841:844").Insert()

# Row 841 - Tomate / Larga vida / Primera (new weekly price entry)
$ws.Range("A841").Value = 2
$ws.Range("B841").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C841").Value = "Coquimbo"
$ws.Range("D841").Value = 45021
$ws.Range("E841").Value = 4
$ws.Range("F841").Value = 100112020
$ws.Range("G841").Value = "Tomate"
$ws.Range("H841").Value = "Larga vida"
$ws.Range("I841").Value = "Primera"
$ws.Range("J841").Value = 2400
$ws.Range("K841").Value = 5000
$ws.Range("L841").Value = 6000
$ws.Range("M841").Value = 5500
$ws.Range("N841").Value = "$/bandeja 18 kilos"
$ws.Range("O841").Value = "Provincia de Limarí"
$ws.Range("P841").Value = 306
$ws.Range("Q841").Value = 18
$ws.Range("R841").Value = "Hortaliza"

# Row 842 - Tomate / Larga vida / Segunda
$ws.Range("A842").Value = 2
$ws.Range("B842").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C842").Value = "Coquimbo"
$ws.Range("D842").Value = 45021
$ws.Range("E842").Value = 4
$ws.Range("F842").Value = 100112020
$ws.Range("G842").Value = "Tomate"
$ws.Range("H842").Value = "Larga vida"
$ws.Range("I842").Value = "Segunda"
$ws.Range("J842").Value = 2000
$ws.Range("K842").Value = 3000
$ws.Range("L842").Value = 4000
$ws.Range("M842").Value = 3500
$ws.Range("N842").Value = "$/bandeja 18 kilos"
$ws.Range("O842").Value = "Provincia de Limarí"
$ws.Range("P842").Value = 194
$ws.Range("Q842").Value = 18
$ws.Range("R842").Value = "Hortaliza"

# Row 843 - Tomate / Semiduro / Primera
$ws.Range("A843").Value = 2
$ws.Range("B843").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C843").Value = "Coquimbo"
$ws.Range("D843").Value = 45021
$ws.Range("E843").Value = 4
$ws.Range("F843").Value = 100112020
$ws.Range("G843").Value = "Tomate"
$ws.Range("H843").Value = "Semiduro"
$ws.Range("I843").Value = "Primera"
$ws.Range("J843").Value = 2800
$ws.Range("K843").Value = 4500
$ws.Range("L843").Value = 5000
$ws.Range("M843").Value = 4750
$ws.Range("N843").Value = "$/bandeja 18 kilos"
$ws.Range("O843").Value = "Provincia de Limarí"
$ws.Range("P843").Value = 264
$ws.Range("Q843").Value = 18
$ws.Range("R843").Value = "Hortaliza"

# Row 844 - Tomate / Semiduro / Segunda
$ws.Range("A844").Value = 2
$ws.Range("B844").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C844").Value = "Coquimbo"
$ws.Range("D844").Value = 45021
$ws.Range("E844").Value = 4
$ws.Range("F844").Value = 100112020
$ws.Range("G844").Value = "Tomate"
$ws.Range("H844").Value = "Semiduro"
$ws.Range("I844").Value = "Segunda"
$ws.Range("J844").Value = 2200
$ws.Range("K844").Value = 2500
$ws.Range("L844").Value = 3000
$ws.Range("M844").Value = 2750
$ws.Range("N844").Value = "$/bandeja 18 kilos"
$ws.Range("O844").Value = "Provincia de Limarí"
$ws.Range("P844").Value = 153
$ws.Range("Q844").Value = 18
$ws.Range("R844").Value = "Hortaliza"

# Ensure the date column keeps its date number format on the new rows
# (Excel normally carries formatting down on Insert, but set it explicitly
# too so it is not lost).
$ws.Range("D841:D844").NumberFormat = $ws.Range("D845").NumberFormat
